$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = '(''Akki Coalflinger'', [''{1}{R}{R}'', ''Creature — Goblin Shaman'', ''First strike'', ''{R}, {T}: Attacking creatures gain first strike until end of turn.'', ''2/2''])'
$ws.Range("A3").Value2 = '(''Allosaurus Rider'', [''{5}{G}{G}'', ''Creature — Elf Warrior'', ''You may exile two green cards from your hand rather than pay this spell’s mana cost.'', ''Allosaurus Rider’s power and toughness are each equal to 1 plus the number of lands you control.'', ''1+*/1+*''])'
$ws.Range("A4").Value2 = '(''Ambush Commander'', [''{3}{G}{G}'', ''Creature — Elf'', ''Forests you control are 1/1 green Elf creatures that are still lands.'', ''{1}{G}, Sacrifice an Elf: Target creature gets +3/+3 until end of turn.'', ''2/2''])'
$ws.Range("A5").Value2 = '(''Boggart Shenanigans'', [''{2}{R}'', ''Tribal Enchantment — Goblin'', ''Whenever another Goblin you control is put into a graveyard from the battlefield, you may have Boggart Shenanigans deal 1 damage to target player or planeswalker.''])'
$ws.Range("A6").Value2 = '(''Clickslither'', [''{1}{R}{R}{R}'', ''Creature — Insect'', ''Haste'', ''Sacrifice a Goblin: Clickslither gets +2/+2 and gains trample until end of turn.'', ''3/3''])'
$ws.Range("A7").Value2 = '(''Elvish Eulogist'', [''{G}'', ''Creature — Elf Shaman'', ''Sacrifice Elvish Eulogist: You gain 1 life for each Elf card in your graveyard.'', ''1/1''])'
$ws.Range("A8").Value2 = '(''Elvish Harbinger'', [''{2}{G}'', ''Creature — Elf Druid'', ''When Elvish Harbinger enters the battlefield, you may search your library for an Elf card, reveal it, then shuffle your library and put that card on top of it.'', ''{T}: Add one mana of any color.'', ''1/2''])'
$ws.Range("A9").Value2 = '(''Elvish Promenade'', [''{3}{G}'', ''Tribal Sorcery — Elf'', ''Create a 1/1 green Elf Warrior creature token for each Elf you control.''])'
$ws.Range("A10").Value2 = '(''Elvish Warrior'', [''{G}{G}'', ''Creature — Elf Warrior'', ''2/3''])'
$ws.Range("A11").Value2 = '(''Emberwilde Augur'', [''{1}{R}'', ''Creature — Goblin Shaman'', ''Sacrifice Emberwilde Augur: It deals 3 damage to target player or planeswalker. Activate this ability only during your upkeep.'', ''2/1''])'
$ws.Range("A12").Value2 = '(''Flamewave Invoker'', [''{2}{R}'', ''Creature — Goblin Mutant'', ''{7}{R}: Flamewave Invoker deals 5 damage to target player or planeswalker.'', ''2/2''])'
$ws.Range("A13").Value2 = '(''Forest'', [''Basic Land — Forest'', ''({T}: Add {G}.)''])'
$ws.Range("A14").Value2 = '(''Forgotten Cave'', [''Land'', ''Forgotten Cave enters the battlefield tapped.'', ''{T}: Add {R}.'', ''Cycling {R} ({R}, Discard this card: Draw a card.)''])'
$ws.Range("A15").Value2 = '(''Gempalm Incinerator'', [''{2}{R}'', ''Creature — Goblin'', ''Cycling {1}{R} ({1}{R}, Discard this card: Draw a card.)'', ''When you cycle Gempalm Incinerator, you may have it deal X damage to target creature, where X is the number of Goblins on the battlefield.'', ''2/1''])'
$ws.Range("A16").Value2 = '(''Gempalm Strider'', [''{1}{G}'', ''Creature — Elf'', ''Cycling {2}{G}{G} ({2}{G}{G}, Discard this card: Draw a card.)'', ''When you cycle Gempalm Strider, Elf creatures get +2/+2 until end of turn.'', ''2/2''])'
$ws.Range("A17").Value2 = '(''Giant Growth'', [''{G}'', ''Instant'', ''Target creature gets +3/+3 until end of turn.''])'
$ws.Range("A18").Value2 = '(''Goblin Burrows'', [''Land'', ''{T}: Add {C}.'', ''{1}{R}, {T}: Target Goblin creature gets +2/+0 until end of turn.''])'
$ws.Range("A19").Value2 = '(''Goblin Cohort'', [''{R}'', ''Creature — Goblin Warrior'', ''Goblin Cohort can’t attack unless you’ve cast a creature spell this turn.'', ''2/2''])'
$ws.Range("A20").Value2 = '(''Goblin Matron'', [''{2}{R}'', ''Creature — Goblin'', ''When Goblin Matron enters the battlefield, you may search your library for a Goblin card, reveal that card, put it into your hand, then shuffle your library.'', ''1/1''])'
$ws.Range("A21").Value2 = '(''Goblin Ringleader'', [''{3}{R}'', ''Creature — Goblin'', ''Haste (This creature can attack and {T} as soon as it comes under your control.)'', ''When Goblin Ringleader enters the battlefield, reveal the top four cards of your library. Put all Goblin cards revealed this way into your hand and the rest on the bottom of your library in any order.'', ''2/2''])'
$ws.Range("A22").Value2 = '(''Goblin Sledder'', [''{R}'', ''Creature — Goblin'', ''Sacrifice a Goblin: Target creature gets +1/+1 until end of turn.'', ''1/1''])'
$ws.Range("A23").Value2 = '(''Goblin Warchief'', [''{1}{R}{R}'', ''Creature — Goblin Warrior'', ''Goblin spells you cast cost {1} less to cast.'', ''Goblins you control have haste.'', ''2/2''])'
$ws.Range("A24").Value2 = '(''Harmonize'', [''{2}{G}{G}'', ''Sorcery'', ''Draw three cards.''])'
$ws.Range("A25").Value2 = '(''Heedless One'', [''{3}{G}'', ''Creature — Elf Avatar'', ''Trample'', ''Heedless One’s power and toughness are each equal to the number of Elves on the battlefield.'', ''*/*''])'
$ws.Range("A26").Value2 = '(''Ib Halfheart, Goblin Tactician'', [''{3}{R}'', ''Legendary Creature — Goblin Advisor'', ''Whenever another Goblin you control becomes blocked, sacrifice it. If you do, it deals 4 damage to each creature blocking it.'', ''Sacrifice two Mountains: Create two 1/1 red Goblin creature tokens.'', ''3/2''])'
$ws.Range("A27").Value2 = '(''Imperious Perfect'', [''{2}{G}'', ''Creature — Elf Warrior'', ''Other Elves you control get +1/+1.'', ''{G}, {T}: Create a 1/1 green Elf Warrior creature token.'', ''2/2''])'
$ws.Range("A28").Value2 = '(''Llanowar Elves'', [''{G}'', ''Creature — Elf Druid'', ''{T}: Add {G}.'', ''1/1''])'
$ws.Range("A29").Value2 = '(''Lys Alana Huntmaster'', [''{2}{G}{G}'', ''Creature — Elf Warrior'', ''Whenever you cast an Elf spell, you may create a 1/1 green Elf Warrior creature token.'', ''3/3''])'
$ws.Range("A30").Value2 = '(''Mogg Fanatic'', [''{R}'', ''Creature — Goblin'', ''Sacrifice Mogg Fanatic: It deals 1 damage to any target.'', ''1/1''])'
$ws.Range("A31").Value2 = '(''Mogg War Marshal'', [''{1}{R}'', ''Creature — Goblin Warrior'', ''Echo {1}{R} (At the beginning of your upkeep, if this came under your control since the beginning of your last upkeep, sacrifice it unless you pay its echo cost.)'', ''When Mogg War Marshal enters the battlefield or dies, create a 1/1 red Goblin creature token.'', ''1/1''])'
$ws.Range("A32").Value2 = '(''Moonglove Extract'', [''{3}'', ''Artifact'', ''Sacrifice Moonglove Extract: It deals 2 damage to any target.''])'
$ws.Range("A33").Value2 = '(''Mountain'', [''Basic Land — Mountain'', ''({T}: Add {R}.)''])'
$ws.Range("A34").Value2 = '(''Mudbutton Torchrunner'', [''{2}{R}'', ''Creature — Goblin Warrior'', ''When Mudbutton Torchrunner dies, it deals 3 damage to any target.'', ''1/1''])'
$ws.Range("A35").Value2 = '(''Raging Goblin'', [''{R}'', ''Creature — Goblin Berserker'', ''Haste (This creature can attack and {T} as soon as it comes under your control.)'', ''1/1''])'
$ws.Range("A36").Value2 = '(''Reckless One'', [''{3}{R}'', ''Creature — Goblin Avatar'', ''Haste'', ''Reckless One’s power and toughness are each equal to the number of Goblins on the battlefield.'', ''*/*''])'
$ws.Range("A37").Value2 = '(''Siege-Gang Commander'', [''{3}{R}{R}'', ''Creature — Goblin'', ''When Siege-Gang Commander enters the battlefield, create three 1/1 red Goblin creature tokens.'', ''{1}{R}, Sacrifice a Goblin: Siege-Gang Commander deals 2 damage to any target.'', ''2/2''])'
$ws.Range("A38").Value2 = '(''Skirk Drill Sergeant'', [''{1}{R}'', ''Creature — Goblin'', ''Whenever Skirk Drill Sergeant or another Goblin dies, you may pay {2}{R}. If you do, reveal the top card of your library. If it’s a Goblin permanent card, put it onto the battlefield. Otherwise, put it into your graveyard.'', ''2/1''])'
$ws.Range("A39").Value2 = '(''Skirk Fire Marshal'', [''{3}{R}{R}'', ''Creature — Goblin'', ''Protection from red'', ''Tap five untapped Goblins you control: Skirk Fire Marshal deals 10 damage to each creature and each player.'', ''2/2''])'
$ws.Range("A40").Value2 = '(''Skirk Prospector'', [''{R}'', ''Creature — Goblin'', ''Sacrifice a Goblin: Add {R}.'', ''1/1''])'
$ws.Range("A41").Value2 = '(''Skirk Shaman'', [''{1}{R}{R}'', ''Creature — Goblin Shaman'', ''Skirk Shaman can’t be blocked except by artifact creatures and/or red creatures.'', ''2/2''])'
$ws.Range("A42").Value2 = '(''Slate of Ancestry'', [''{4}'', ''Artifact'', ''{4}, {T}, Discard your hand: Draw a card for each creature you control.''])'
$ws.Range("A43").Value2 = '(''Spitting Earth'', [''{1}{R}'', ''Sorcery'', ''Spitting Earth deals damage to target creature equal to the number of Mountains you control.''])'
$ws.Range("A44").Value2 = '(''Stonewood Invoker'', [''{1}{G}'', ''Creature — Elf Mutant'', ''{7}{G}: Stonewood Invoker gets +5/+5 until end of turn.'', ''2/2''])'
$ws.Range("A45").Value2 = '(''Sylvan Messenger'', [''{3}{G}'', ''Creature — Elf'', ''Trample (This creature can deal excess combat damage to the player or planeswalker it’s attacking.)'', ''When Sylvan Messenger enters the battlefield, reveal the top four cards of your library. Put all Elf cards revealed this way into your hand and the rest on the bottom of your library in any order.'', ''2/2''])'
$ws.Range("A46").Value2 = '(''Tarfire'', [''{R}'', ''Tribal Instant — Goblin'', ''Tarfire deals 2 damage to any target.''])'
$ws.Range("A47").Value2 = '(''Tar Pitcher'', [''{3}{R}'', ''Creature — Goblin Shaman'', ''{T}, Sacrifice a Goblin: Tar Pitcher deals 2 damage to any target.'', ''2/2''])'
$ws.Range("A48").Value2 = '(''Timberwatch Elf'', [''{2}{G}'', ''Creature — Elf'', ''{T}: Target creature gets +X/+X until end of turn, where X is the number of Elves on the battlefield.'', ''1/2''])'
$ws.Range("A49").Value2 = '(''Tranquil Thicket'', [''Land'', ''Tranquil Thicket enters the battlefield tapped.'', ''{T}: Add {G}.'', ''Cycling {G} ({G}, Discard this card: Draw a card.)''])'
$ws.Range("A50").Value2 = '(''Voice of the Woods'', [''{3}{G}{G}'', ''Creature — Elf'', ''Tap five untapped Elves you control: Create a 7/7 green Elemental creature token with trample.'', ''2/2''])'
$ws.Range("A51").Value2 = '(''Wellwisher'', [''{1}{G}'', ''Creature — Elf'', ''{T}: You gain 1 life for each Elf on the battlefield.'', ''1/1''])'
$ws.Range("A52").Value2 = '(''Wildsize'', [''{2}{G}'', ''Instant'', ''Target creature gets +2/+2 and gains trample until end of turn.'', ''Draw a card.''])'
$ws.Range("A53").Value2 = '(''Wirewood Herald'', [''{1}{G}'', ''Creature — Elf'', ''When Wirewood Herald dies, you may search your library for an Elf card, reveal that card, put it into your hand, then shuffle your library.'', ''1/1''])'
$ws.Range("A54").Value2 = '(''Wirewood Lodge'', [''Land'', ''{T}: Add {C}.'', ''{G}, {T}: Untap target Elf.''])'
$ws.Range("A55").Value2 = '(''Wirewood Symbiote'', [''{G}'', ''Creature — Insect'', ''Return an Elf you control to its owner’s hand: Untap target creature. Activate this ability only once each turn.'', ''1/1''])'
$ws.Range("A56").Value2 = '(''Wood Elves'', [''{2}{G}'', ''Creature — Elf Scout'', ''When Wood Elves enters the battlefield, search your library for a Forest card and put that card onto the battlefield. Then shuffle your library.'', ''1/1''])'
$ws.Range("A57").Value2 = '("Wren''s Run Vanquisher", [''{1}{G}'', ''Creature — Elf Warrior'', ''As an additional cost to cast this spell, reveal an Elf card from your hand or pay {3}.'', ''Deathtouch (Any amount of damage this deals to a creature is enough to destroy it.)'', ''3/3''])'
